$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted before the existing row 422,
# pushing the former rows 422-512 down to 423-513 (dimension A1:R512 -> A1:R513).
$ws.Rows.Item(422).Insert()

# Populate the newly inserted row 422 with the new record. It mirrors the
# record that used to be at row 422 (same market/category/price/etc.) except
# for a newer date and an updated Origen value.
$ws.Cells.Item(422, 1).Value = 9
$ws.Cells.Item(422, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(422, 3).Value = "Metropolitana"
$ws.Cells.Item(422, 4).Value = 45275
$ws.Cells.Item(422, 5).Value = 13
$ws.Cells.Item(422, 6).Value = 300000001
$ws.Cells.Item(422, 7).Value = "Rabanito"
$ws.Cells.Item(422, 8).Value = "Sin especificar"
$ws.Cells.Item(422, 9).Value = "Primera"
$ws.Cells.Item(422, 10).Value = 7000
$ws.Cells.Item(422, 11).Value = 3000
$ws.Cells.Item(422, 12).Value = 3000
$ws.Cells.Item(422, 13).Value = 3000
$ws.Cells.Item(422, 14).Value = "$/cien unidades (volumen en unidades)"
$ws.Cells.Item(422, 15).Value = "Región Metropolitana"
$ws.Cells.Item(422, 16).Value = 30
$ws.Cells.Item(422, 17).Value = 100
$ws.Cells.Item(422, 18).Value = "Hortaliza"
